# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a number of rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 8;   Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 10;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 56;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 59;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 66;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 67;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 76;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 88;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 95;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 100; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 101; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 109; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 111; Tag = "b";  Act = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
